$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K" = strikeouts),
# replacing the old "Strike#" derived values with the correct K counts.
$kValues = @{
    2 = 2
    3 = 0
    4 = 1
    5 = 0
    6 = 1
    7 = 0
    8 = 1
    9 = 1
    10 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    24 = 0
    25 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 0
    31 = 1
    33 = 1
    34 = 2
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 0
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    45 = 1
    49 = 1
    50 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
